{"js": "// The source document has a run-of-text (\"FE-10: A staff can view reports to\n// enable or disable accounts of landlords\") that ends with two separate\n// runs sharing identical formatting:\n//    <w:r> <w:t xml:space=\"preserve\"> </w:t> </w:r>\n//    <w:r> <w:t>enable or disable accounts of landlords</w:t> </w:r>\n// The edit folds those two runs into a single run:\n//    <w:r> <w:t xml:space=\"preserve\"> enable or disable accounts of landlords</w:t> </w:r>\n// (the visible text is unchanged - only the run split disappears).\n//\n// A plain Range.insertText()/Range.Text= replace would also touch the\n// *preceding* sibling run (\" view reports to\") because every run in this\n// paragraph shares the same rPr and the host coalesces adjacent\n// same-format runs around any edit point. To keep that earlier run intact\n// (matching the target XML exactly) we splice in literal OOXML for just\n// the two-run span that must become one run, using Range.insertOoxml,\n// which replaces only the addressed range's contents verbatim.\n\nconst searchText = \" enable or disable accounts of landlords\";\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + searchText);\n}\n\nconst target = results.items[0];\n\nconst mergedRunOoxml =\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>' +\n  '<w:sz w:val=\"22\"/>' +\n  '<w:szCs w:val=\"22\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\"> enable or disable accounts of landlords</w:t>' +\n  '</w:r>';\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' + mergedRunOoxml + '</w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(flatOpcPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\nWrite-Output $d.Content.Text.Length\n"}
